$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 133.7780026666667
$ws.Range("N2").Value = 401.334008
$ws.Range("O2").Value = 0.50863533211804
$ws.Range("P2").Value = 0.5086353321180399
$ws.Range("Q2").Value = 41.89646109714401
$ws.Range("R2").Value = 377.068149874296
$ws.Range("S2").Value = 0.01377112154277121
$ws.Range("T2").Value = 0.01377112154277121

$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.1993888292903622
$ws.Range("P3").Value = 0.1993888292903622
$ws.Range("Q3").Value = 16.423724035809
$ws.Range("R3").Value = 147.813516322281
$ws.Range("S3").Value = 0.005398381962563335
$ws.Range("T3").Value = 0.005398381962563335

$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 21.197691
$ws.Range("N4").Value = 63.593073
$ws.Range("O4").Value = 0.08059542216956049
$ws.Range("P4").Value = 0.08059542216956046
$ws.Range("Q4").Value = 6.638671669689001
$ws.Range("R4").Value = 59.748045027201
$ws.Range("S4").Value = 0.002182092521701581
$ws.Range("T4").Value = 0.00218209252170158

$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 55.59592133333333
$ws.Range("N5").Value = 166.787764
$ws.Range("O5").Value = 0.2113804164220374
$ws.Range("P5").Value = 0.2113804164220373
$ws.Range("Q5").Value = 17.411475047252
$ws.Range("R5").Value = 156.703275425268
$ws.Range("S5").Value = 0.005723049938720968
$ws.Range("T5").Value = 0.005723049938720967

$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 133.7780026666667
$ws.Range("N6").Value = 401.334008
$ws.Range("O6").Value = 0.50863533211804
$ws.Range("P6").Value = 0.5086353321180399
$ws.Range("Q6").Value = 1080.549185949817
$ws.Range("R6").Value = 9724.942673548352
$ws.Range("S6").Value = 0.3551701929705892
$ws.Range("T6").Value = 0.3551701929705892

$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.1993888292903622
$ws.Range("P7").Value = 0.1993888292903622
$ws.Range("S7").Value = 0.1392293545168097
$ws.Range("T7").Value = 0.1392293545168097

$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 21.197691
$ws.Range("N8").Value = 63.593073
$ws.Range("O8").Value = 0.08059542216956049
$ws.Range("P8").Value = 0.08059542216956046
$ws.Range("Q8").Value = 171.217594054968
$ws.Range("R8").Value = 1540.958346494712
$ws.Range("S8").Value = 0.0562782210298081
$ws.Range("T8").Value = 0.05627822102980808

$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 55.59592133333333
$ws.Range("N9").Value = 166.787764
$ws.Range("O9").Value = 0.2113804164220374
$ws.Range("P9").Value = 0.2113804164220373
$ws.Range("Q9").Value = 449.0583380030685
$ws.Range("R9").Value = 4041.525042027616
$ws.Range("S9").Value = 0.1476028473645151
$ws.Range("T9").Value = 0.1476028473645151

$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 133.7780026666667
$ws.Range("N10").Value = 401.334008
$ws.Range("O10").Value = 0.50863533211804
$ws.Range("P10").Value = 0.5086353321180399
$ws.Range("Q10").Value = 387.6546275226552
$ws.Range("R10").Value = 3488.891647703897
$ws.Range("S10").Value = 0.1274198071253349
$ws.Range("T10").Value = 0.1274198071253349

$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.1993888292903622
$ws.Range("P11").Value = 0.1993888292903622
$ws.Range("Q11").Value = 151.963494216709
$ws.Range("R11").Value = 1367.671447950381
$ws.Range("S11").Value = 0.04994951110716045
$ws.Range("T11").Value = 0.04994951110716044

$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 21.197691
$ws.Range("N12").Value = 63.593073
$ws.Range("O12").Value = 0.08059542216956049
$ws.Range("P12").Value = 0.08059542216956046
$ws.Range("Q12").Value = 61.42551723858901
$ws.Range("R12").Value = 552.8296551473011
$ws.Range("S12").Value = 0.02019020799295769
$ws.Range("T12").Value = 0.02019020799295768

$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 55.59592133333333
$ws.Range("N13").Value = 166.787764
$ws.Range("O13").Value = 0.2113804164220374
$ws.Range("P13").Value = 0.2113804164220373
$ws.Range("Q13").Value = 161.1028401280076
$ws.Range("R13").Value = 1449.925561152068
$ws.Range("S13").Value = 0.05295356061564032
$ws.Range("T13").Value = 0.05295356061564031

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 133.7780026666667
$ws.Range("N14").Value = 401.334008
$ws.Range("O14").Value = 0.50863533211804
$ws.Range("P14").Value = 0.5086353321180399
$ws.Range("Q14").Value = 37.34234573769778
$ws.Range("R14").Value = 336.08111163928
$ws.Range("S14").Value = 0.01227421047934465
$ws.Range("T14").Value = 0.01227421047934465

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.1993888292903622
$ws.Range("P15").Value = 0.1993888292903622
$ws.Range("Q15").Value = 14.63847697837
$ws.Range("R15").Value = 131.74629280533
$ws.Range("S15").Value = 0.004811581703828761
$ws.Range("T15").Value = 0.004811581703828761

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 21.197691
$ws.Range("N16").Value = 63.593073
$ws.Range("O16").Value = 0.08059542216956049
$ws.Range("P16").Value = 0.08059542216956046
$ws.Range("Q16").Value = 5.91705280677
$ws.Range("R16").Value = 53.25347526093
$ws.Range("S16").Value = 0.001944900625093126
$ws.Range("T16").Value = 0.001944900625093126

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 55.59592133333333
$ws.Range("N17").Value = 166.787764
$ws.Range("O17").Value = 0.2113804164220374
$ws.Range("P17").Value = 0.2113804164220373
$ws.Range("Q17").Value = 15.51886016124889
$ws.Range("R17").Value = 139.66974145124
$ws.Range("S17").Value = 0.005100958503160946
$ws.Range("T17").Value = 0.005100958503160946
